# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullets (the ones that follow
# the "Impact" Heading3, immediately before "TECHNICAL SKILLS") as
# impact-focused accomplishment statements, and drop the FEC-analysis bullet
# entirely (5 bullets remain instead of 6).
#
# There is a near-duplicate set of bullets earlier in the document (under
# "PROFESSIONAL EXPERIENCE" / "Partner - Siege Analytics") that must be left
# untouched, so we locate paragraphs by walking from the "KEY ACHIEVEMENTS
# AND IMPACT" heading rather than doing a blind Find/Replace across the
# whole document.

$d = $word.ActiveDocument
$paras = $d.Paragraphs
$count = $paras.Count

# Find the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph.
$headingIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($paras.Item($i).Range.Text.Trim() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# The very next paragraph is the "Impact" Heading3 sub-label; the bullet
# paragraphs start right after that.
$firstBulletIndex = $headingIndex + 2

# Bullet 1: Discovered systematic race coding errors ... -> Breakthrough demographic discovery
$p1 = $paras.Item($firstBulletIndex)
if ($p1.Range.Text.Trim() -ne "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%") {
    throw "Unexpected text at bullet 1: $($p1.Range.Text)"
}
$p1.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# Bullet 2: Trigonometric algorithm for boundary estimation ... -> 178% accuracy improvement
$p2 = $paras.Item($firstBulletIndex + 1)
if ($p2.Range.Text.Trim() -ne "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis") {
    throw "Unexpected text at bullet 2: $($p2.Range.Text)"
}
$p2.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Bullet 3: Built redistricting platform ... -> Algorithmic innovation: trigonometric boundary estimation
$p3 = $paras.Item($firstBulletIndex + 2)
if ($p3.Range.Text.Trim() -ne "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations") {
    throw "Unexpected text at bullet 3: $($p3.Range.Text)"
}
$p3.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"

# Bullet 4: Achieved 87% prediction accuracy ... -> $4.7M savings enabled nonprofit access
$p4 = $paras.Item($firstBulletIndex + 3)
if ($p4.Range.Text.Trim() -ne "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%") {
    throw "Unexpected text at bullet 4: $($p4.Range.Text)"
}
$p4.Range.Text = "• `$4.7M savings enabled nonprofit access"

# Bullet 5: Built real-time FEC analysis systems ... -> removed entirely
$p5 = $paras.Item($firstBulletIndex + 4)
if ($p5.Range.Text.Trim() -ne "• Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion") {
    throw "Unexpected text at bullet 5: $($p5.Range.Text)"
}
$p5.Range.Delete()

# Bullet 6 (now shifted up to firstBulletIndex + 4 after the delete):
# Provided expert testimony ... -> Platform impact: Built redistricting system
$p6 = $paras.Item($firstBulletIndex + 4)
if ($p6.Range.Text.Trim() -ne "• Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy") {
    throw "Unexpected text at bullet 6: $($p6.Range.Text)"
}
$p6.Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

Write-Host "Key Achievements section updated successfully"
